# Apply "Final analysis and results" edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header renames
$ws.Range("B1").Value = "sparsity_necessary"
$ws.Range("C1").Value = "necessary explanation rate"

# Rows 2-11: selection strategy "1-best" -> "1-delta"
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("E$r").Value = "1-delta"
}

# Rows 39-45: selection strategy "1-best" -> "temporal"
for ($r = 39; $r -le 45; $r++) {
    $ws.Range("E$r").Value = "temporal"
}

# Updated numeric values in the "Greedy" block
$ws.Range("C40").Value = 0.175
$ws.Range("C41").Value = 0.27
$ws.Range("C42").Value = 0.31
$ws.Range("C43").Value = 0.32
$ws.Range("B44").Value = 0.09375
$ws.Range("C44").Value = 0.33
$ws.Range("C45").Value = 0.33
